# NV128.xlsx update
# - Add four new "RIVA128 to do" list items (Chroma Key, Pattern,
#   Unify position/size structs, Code cleanup) as new rows below the
#   existing list.
# - Leave the cursor/selection on the cell the author last had selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A35").Value = "Chroma Key"
$ws.Range("A36").Value = "Pattern"
$ws.Range("A37").Value = "Unify position/size structs"
$ws.Range("A38").Value = "Code cleanup"

# Match the author's final selection/scroll position.
$ws.Range("E31").Select()
